# Weekly fruit/vegetable data update: a new price-report row for
# "Vega Modelo de Temuco - Haba" is inserted at the top of the data
# (row 65, just under the header), pushing the previously-existing
# rows 65-75 down to rows 66-76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data row 65 (and everything below it) down by one row.
$ws.Rows("65").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A65").Value = 10
$ws.Range("B65").Value = "Vega Modelo de Temuco"
$ws.Range("C65").Value = "La Araucanía"
$ws.Range("D65").Value = 44826
$ws.Range("E65").Value = 9
$ws.Range("F65").Value = 100112026
$ws.Range("G65").Value = "Haba"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 180
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = 12778
$ws.Range("N65").Value = '$/saco 25 kilos'
$ws.Range("O65").Value = "Región Metropolitana"
$ws.Range("P65").Value = 511
$ws.Range("Q65").Value = 25
$ws.Range("R65").Value = "Hortaliza"
